# Update countries & provincias Spain
# Refreshes the COVID-19 country snapshot on sheet "Pais": new totals for a
# handful of countries (time-series refresh) plus three countries whose
# growth pushed them past their neighbours in the (descending, by total
# cases) ranking, so their rows now carry different countries/values than
# before. Finally the "datos actualizados" timestamp is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Header timestamp: 09:34 -> 10:51
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 10:51"

# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

# Rusia (row 6)
$ws.Range("B6").Value = 599705
$ws.Range("C6").Value = 7425
$ws.Range("D6").Value = 356429
$ws.Range("E6").Value = 234917
$ws.Range("G6").Value = 153
$ws.Range("H6").Value = 8359

# Banglades (row 20)
$ws.Range("B20").Value = 119198
$ws.Range("C20").Value = 3412
$ws.Range("D20").Value = 47635
$ws.Range("E20").Value = 70018
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = 1545

# Polonia (row 40)
$ws.Range("B40").Value = 32527
$ws.Range("C40").Value = 300
$ws.Range("D40").Value = 17573
$ws.Range("E40").Value = 13579
$ws.Range("G40").Value = 16
$ws.Range("H40").Value = 1375

# Barein (row 50)
$ws.Range("E50").Value = 5479
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 66

# Israel (row 51)
$ws.Range("B51").Value = 21246
$ws.Range("C51").Value = 164
$ws.Range("D51").Value = 15812
$ws.Range("E51").Value = 5127

# Austria (row 56)
$ws.Range("B56").Value = 17408
$ws.Range("C56").Value = 28
$ws.Range("D56").Value = 16261
$ws.Range("E56").Value = 454
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 693

# El Salvador (row 83)
$ws.Range("B83").Value = 4973
$ws.Range("C83").Value = 165
$ws.Range("D83").Value = 2814
$ws.Range("E83").Value = 2046
$ws.Range("G83").Value = 6
$ws.Range("H83").Value = 113

# Eslovaquia (row 116)
$ws.Range("B116").Value = 1589
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 1448

# Estado de Palestina overtakes Niger and Jordania (rows 127-129):
# row 127 now shows Estado de Palestina's refreshed numbers, and Niger /
# Jordania shift down one row each, keeping their own totals unchanged.
$ws.Range("A127").Value = "Estado de Palestina"
$ws.Range("B127").Value = 1078
$ws.Range("C127").Value = 77
$ws.Range("D127").Value = 442
$ws.Range("E127").Value = 633
$ws.Range("H127").Value = 3

$ws.Range("A128").Value = "Niger"
$ws.Range("B128").Value = 1046
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 913
$ws.Range("E128").Value = 66
$ws.Range("H128").Value = 67

$ws.Range("A129").Value = "Jordania"
$ws.Range("B129").Value = 1042
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 751
$ws.Range("E129").Value = 282
$ws.Range("H129").Value = 9

# Malaui overtakes Uganda and Ruanda (rows 139-141): row 139 now shows
# Malaui's refreshed numbers, Uganda keeps its row with refreshed numbers
# too, and Ruanda shifts down to row 141 keeping its own totals unchanged.
$ws.Range("A139").Value = "Malaui"
$ws.Range("B139").Value = 803
$ws.Range("C139").Value = 54
$ws.Range("D139").Value = 258
$ws.Range("E139").Value = 534
$ws.Range("H139").Value = 11

$ws.Range("A140").Value = "Uganda"
$ws.Range("B140").Value = 797
$ws.Range("C140").Value = 23
$ws.Range("D140").Value = 699
$ws.Range("E140").Value = 98

$ws.Range("A141").Value = "Ruanda"
$ws.Range("B141").Value = 787
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 370
$ws.Range("E141").Value = 415
$ws.Range("H141").Value = 2
